# Update gh-pages to output generated at 456a3b4
# "想去人数" (want-to-go count) and related figures refreshed for the
# first event (南宁·2024良牙动漫冬季盛典) and the third event
# (南宁·0316全职only-全明星周末) on both the "展览" sheet and the
# aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 9627
    $ws.Range("F4").Value = 32
}
